$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: "Save" (reuse same format/style as the other header cells, e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data column H2:H6, all zeros
$ws.Range("H2:H6").Value = 0
